# Updated code and protocol
# Refresh the letter-frequency table (column A labels + B counts + C shares)
# with a new corpus pass, adding one more row (35) for the newly observed letter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35 is brand new -- clone row 34s formatting (border/bold/centered) onto it
# before we repurpose row 34 for different data, so the new row keeps the same style.
$ws.Range("A34").Copy($ws.Range("A35"))

$ws.Range("A2").Value = " "
$ws.Range("B2").Value = 205592
$ws.Range("C2").Value = 0.166714104189017
$ws.Range("A3").Value = "о"
$ws.Range("B3").Value = 115667
$ws.Range("C3").Value = 0.09379411790940811
$ws.Range("A4").Value = "е"
$ws.Range("B4").Value = 91838
$ws.Range("C4").Value = 0.07447123380535696
$ws.Range("A5").Value = "а"
$ws.Range("B5").Value = 81731
$ws.Range("C5").Value = 0.06627548955928514
$ws.Range("A6").Value = "н"
$ws.Range("B6").Value = 68872
$ws.Range("C6").Value = 0.05584815451820101
$ws.Range("A7").Value = "и"
$ws.Range("B7").Value = 65574
$ws.Range("C7").Value = 0.05317381351458522
$ws.Range("A8").Value = "т"
$ws.Range("B8").Value = 64769
$ws.Range("C8").Value = 0.05252104077113139
$ws.Range("A9").Value = "с"
$ws.Range("B9").Value = 54653
$ws.Range("C9").Value = 0.04431799844469798
$ws.Range("A10").Value = "в"
$ws.Range("B10").Value = 48418
$ws.Range("C10").Value = 0.03926205054974818
$ws.Range("A11").Value = "л"
$ws.Range("B11").Value = 47693
$ws.Range("C11").Value = 0.03867414963173076
$ws.Range("A12").Value = "р"
$ws.Range("B12").Value = 40108
$ws.Range("C12").Value = 0.03252348968254162
$ws.Range("A13").Value = "к"
$ws.Range("B13").Value = 33391
$ws.Range("C13").Value = 0.02707668903933746
$ws.Range("A14").Value = "м"
$ws.Range("B14").Value = 32127
$ws.Range("C14").Value = 0.02605171419744227
$ws.Range("A15").Value = "д"
$ws.Range("B15").Value = 31566
$ws.Range("C15").Value = 0.02559680052156948
$ws.Range("A16").Value = "п"
$ws.Range("B16").Value = 27828
$ws.Range("C16").Value = 0.02256566447805346
$ws.Range("A17").Value = "у"
$ws.Range("B17").Value = 26970
$ws.Range("C17").Value = 0.02186991415024801
$ws.Range("A18").Value = "я"
$ws.Range("B18").Value = 24528
$ws.Range("C18").Value = 0.01988970167880175
$ws.Range("A19").Value = "ь"
$ws.Range("B19").Value = 23157
$ws.Range("C19").Value = 0.01877796077038536
$ws.Range("A20").Value = "ч"
$ws.Range("B20").Value = 19631
$ws.Range("C20").Value = 0.01591873506427582
$ws.Range("A21").Value = "г"
$ws.Range("B21").Value = 19295
$ws.Range("C21").Value = 0.01564627339744291
$ws.Range("A22").Value = "з"
$ws.Range("B22").Value = 17831
$ws.Range("C22").Value = 0.0144591189919567
$ws.Range("A23").Value = "ы"
$ws.Range("B23").Value = 17825
$ws.Range("C23").Value = 0.01445425360504897
$ws.Range("A24").Value = "б"
$ws.Range("B24").Value = 17200
$ws.Range("C24").Value = 0.01394744246882706
$ws.Range("A25").Value = "ж"
$ws.Range("B25").Value = 12116
$ws.Range("C25").Value = 0.009824837962343527
$ws.Range("A26").Value = "й"
$ws.Range("B26").Value = 10105
$ws.Range("C26").Value = 0.008194122450435899
$ws.Range("A27").Value = "ш"
$ws.Range("B27").Value = 8395
$ws.Range("C27").Value = 0.006807487181732743
$ws.Range("A28").Value = "х"
$ws.Range("B28").Value = 7477
$ws.Range("C28").Value = 0.006063082984849996
$ws.Range("A29").Value = "ю"
$ws.Range("B29").Value = 6147
$ws.Range("C29").Value = 0.004984588886969764
$ws.Range("A30").Value = "э"
$ws.Range("B30").Value = 3674
$ws.Range("C30").Value = 0.002979238583166897
$ws.Range("A31").Value = "щ"
$ws.Range("B31").Value = 3024
$ws.Range("C31").Value = 0.002452155001496106
$ws.Range("A32").Value = "ц"
$ws.Range("B32").Value = 2979
$ws.Range("C32").Value = 0.002415664599688129
$ws.Range("A33").Value = "ф"
$ws.Range("B33").Value = 1877
$ws.Range("C33").Value = 0.001522055204301651
$ws.Range("A34").Value = "ё"
$ws.Range("B34").Value = 836
$ws.Range("C34").Value = 0.0006779105758104316
$ws.Range("A35").Value = "ъ"
$ws.Range("B35").Value = 307
$ws.Range("C35").Value = 0.0002489456301122039

Write-Output "done"
